$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $newText) {
    $cell = $table.Cell($rowIndex, 1)
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $newText
}

# Simple single-value replacements
Set-CellText $t 1 "0M"
Set-CellText $t 2 "0M"
Set-CellText $t 3 "0M"
Set-CellText $t 4 "202"
Set-CellText $t 5 "0.00002"
Set-CellText $t 6 "0.00009"
Set-CellText $t 9 "0.00004"
Set-CellText $t 12 "0.00758"

# Collapse the multi-run tab-separated rows down to a single value
Set-CellText $t 44 "100"
Set-CellText $t 45 "0.01"
Set-CellText $t 46 "320"
